# IST price update 2025-12-22 10:52
# A new price-check snapshot (timestamped 2025-12-22 16:18) was scraped and
# inserted as the new "most recent" column, pushing every existing
# timestamp column one column to the right (B->C, C->D, ... X->Y).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; this shifts B:X -> C:Y and keeps
# per-cell formatting (header row style, column width) intact for the
# shifted columns. It also grows the sheet dimension to A1:Y26 automatically.
$ws.Columns("B").Insert()

# The inserted column inherits column A's width by default; restore the
# uniform 21-wide formatting used by every other data column. (ColumnWidth
# is expressed in "characters"; the stored sheet width is ColumnWidth +
# 0.8333333333333334, so 20.166666666666668 round-trips to exactly 21.)
$ws.Columns("B").ColumnWidth = 20.166666666666668

# New snapshot timestamp header.
$ws.Range("B1").Value = "2025-12-22 16:18"

# New snapshot prices (one per SKU row). Most rows repeat the most recent
# known price; rows 3 and 14 carry genuinely new values from this scrape.
$ws.Range("B2").Value = 929
$ws.Range("B3").Value = 569
$ws.Range("B4").Value = 299
$ws.Range("B5").Value = 569
$ws.Range("B6").Value = 499
$ws.Range("B7").Value = 569
$ws.Range("B8").Value = 929
$ws.Range("B9").Value = 299
$ws.Range("B10").Value = 299
$ws.Range("B11").Value = 2997
$ws.Range("B12").Value = 569
$ws.Range("B13").Value = 569
$ws.Range("B14").Value = 794
$ws.Range("B15").Value = 499
$ws.Range("B16").Value = 299
$ws.Range("B17").Value = 929
$ws.Range("B18").Value = 499
$ws.Range("B19").Value = 1299
$ws.Range("B20").Value = 929
$ws.Range("B21").Value = 499
$ws.Range("B22").Value = 299
$ws.Range("B23").Value = 1299
$ws.Range("B24").Value = 929
$ws.Range("B25").Value = 929
$ws.Range("B26").Value = 1299
